$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value = 55106.72117531574
$ws.Range("C7").Value = 54686.72117531574
$ws.Range("C8").Value = 53453.519540056266
$ws.Range("C9").Value = 7560.0
$ws.Range("C10").Value = 7140.0
$ws.Range("C11").Value = 14372.170291875165
$ws.Range("C12").Value = 15478.475137969412
$ws.Range("C13").Value = 229.54364550000003
$ws.Range("C14").Value = 39628.24603734633
$ws.Range("C15").Value = 39208.24603734633
$ws.Range("C16").Value = 32068.246037346336
$ws.Range("C17").Value = 31563.168785940274
$ws.Range("C18").Value = 30760.862391846338
$ws.Range("C19").Value = 1077.8400000000001
$ws.Range("C21").Value = 275.5336059060588
$ws.Range("C23").Value = 540412.32721391
$ws.Range("C24").Value = 536293.5342139099
$ws.Range("C25").Value = 524199.9573974926
$ws.Range("C26").Value = 74138.27399999998
$ws.Range("C27").Value = 70019.48099999997
$ws.Range("C28").Value = 140942.84379281756
$ws.Range("C29").Value = 2251.0541911425744
$ws.Range("C30").Value = 388620.3390021423
$ws.Range("C31").Value = 384501.54600214225
$ws.Range("C32").Value = 314482.06500214234
$ws.Range("C33").Value = 309528.94917464105
$ws.Range("C34").Value = 301661.0111749998
$ws.Range("C35").Value = 10569.999635999997
$ws.Range("C37").Value = 2702.0616363586505

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C6").Value = 9258.0
$ws.Range("D6").Value = 21.76271855194696
$ws.Range("C7").Value = 6390.0
$ws.Range("D7").Value = -15.957682917807187
$ws.Range("C8").Value = 6016.0
$ws.Range("D8").Value = -20.876591617140537
$ws.Range("C9").Value = 6546.0
$ws.Range("D9").Value = -13.905945599368676
$ws.Range("C12").Value = 7082.833333333332
$ws.Range("D12").Value = -6.845426472562652

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C8").Value = 6335.0
$ws.Range("D8").Value = 24.97842223373041
$ws.Range("C9").Value = 6115.0
$ws.Range("D9").Value = 20.638208675495097
$ws.Range("C10").Value = 7434.0
$ws.Range("D10").Value = 46.65976178146043
$ws.Range("C11").Value = 5911.0
$ws.Range("D11").Value = 16.613647012404172
$ws.Range("C12").Value = 7919.0
$ws.Range("D12").Value = 56.22795985302464
$ws.Range("C13").Value = 6736.0
$ws.Range("D13").Value = 32.889447855786585
$ws.Range("C14").Value = 5778.5714285714275
$ws.Range("D14").Value = 14.001063915985885

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C8").Value = 269.0
$ws.Range("D8").Value = -64.62068341923337
$ws.Range("C9").Value = 472.0
$ws.Range("D9").Value = -37.92179395493738
$ws.Range("C10").Value = 461.0
$ws.Range("D10").Value = -39.36853180768248

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C8").Value = 441.0
$ws.Range("D8").Value = -41.99896426721904
$ws.Range("C9").Value = 471.5
$ws.Range("D9").Value = -37.98755476642579

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C3").Value = 1108.6666666666665
$ws.Range("D3").Value = -51.39545355323014
$ws.Range("C9").Value = 705.0
$ws.Range("D9").Value = 85.4454883973269
$ws.Range("C10").Value = 377.0
$ws.Range("D10").Value = -0.8326962754720038
$ws.Range("C12").Value = 554.3333333333333
$ws.Range("C16").Value = 705.0
$ws.Range("D16").Value = 85.4454883973269
$ws.Range("C17").Value = 377.0
$ws.Range("D17").Value = -0.8326962754720038
$ws.Range("C19").Value = 554.3333333333333

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C3").Value = 6915.896865902219
$ws.Range("C9").Value = 3457.9484329511106
$ws.Range("C13").Value = 3457.9484329511106

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 1591.0
$ws.Range("D5").Value = -21.530911691450246
$ws.Range("C6").Value = 2204.0
$ws.Range("D6").Value = 8.702621390347996
$ws.Range("C7").Value = 2514.0
$ws.Range("D7").Value = 23.9920100614042
$ws.Range("C8").Value = 2233.0
$ws.Range("D8").Value = 10.132919040220996
$ws.Range("C9").Value = 2135.5
$ws.Range("D9").Value = 5.324159700130713
